$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.116.47'
$ws.Range('E2').Value = '  +2.30%  '
$ws.Range('D3').Value = '2.294.82'
$ws.Range('E3').Value = '  +3.62%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''251.99'
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('D6').Value = '''0.639'
$ws.Range('E6').Value = '  +3.83%  '
$ws.Range('D7').Value = '''74.21'
$ws.Range('E7').Value = '  +8.77%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.646'
$ws.Range('E9').Value = '  +3.66%  '
$ws.Range('D10').Value = '''39.21'
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').Value = '''0.0983'
$ws.Range('E11').Value = '  +4.64%  '
$ws.Range('D12').Value = '''59.18'
$ws.Range('E12').Value = '  -0.34%  '
$ws.Range('D13').Value = '''7.41'
$ws.Range('E13').Value = '  +5.00%  '
$ws.Range('E14').Value = '  +2.24%  '
$ws.Range('D15').Value = '2.645.50'
$ws.Range('E15').Value = '  +3.86%  '
$ws.Range('D16').Value = '''15.37'
$ws.Range('E16').Value = '  +6.25%  '
$ws.Range('D17').Value = '''0.878'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').Value = '2.313.76'
$ws.Range('E18').Value = '  +4.04%  '
$ws.Range('D19').Value = '43.025.47'
$ws.Range('E19').Value = '  +2.37%  '
$ws.Range('D20').Value = '''0.0000100'
$ws.Range('E20').Value = '  +4.12%  '
$ws.Range('D21').Value = '''6.32'
$ws.Range('E21').Value = '  +2.76%  '
$ws.Range('D22').Value = '''72.63'
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').Value = '''235.59'
$ws.Range('E23').Value = '  +1.71%  '
$ws.Range('E24').Value = '  +9.78%  '
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').Value = '''11.59'
$ws.Range('E26').Value = '  +3.51%  '
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('D29').Value = '''3.63'
$ws.Range('E29').Value = '  -1.64%  '
$ws.Range('E30').Value = '  -0.46%  '
$ws.Range('D31').Value = '''167.20'
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('D32').Value = '''21.09'
$ws.Range('E32').Value = '  +3.10%  '
$ws.Range('D33').Value = '''6.39'
$ws.Range('E33').Value = '  +5.98%  '
$ws.Range('E34').Value = '  +5.11%  '
$ws.Range('E35').Value = '  +4.98%  '
$ws.Range('E36').Value = '  +20.96%  '
$ws.Range('E37').Value = '  +3.39%  '
$ws.Range('D38').Value = '''4.62'
$ws.Range('E38').Value = '  +12.91%  '
$ws.Range('D39').Value = '''4.77'
$ws.Range('E39').Value = '  +3.55%  '
$ws.Range('D40').Value = '''0.0309'
$ws.Range('E40').Value = '  -1.44%  '
$ws.Range('E41').Value = '  +20.33%  '
$ws.Range('D42').Value = '''2.35'
$ws.Range('E42').Value = '  +5.21%  '
$ws.Range('D43').Value = '''5.97'
$ws.Range('E43').Value = '  +5.50%  '
$ws.Range('E44').Value = '  +10.36%  '
$ws.Range('D45').Value = '''9.13'
$ws.Range('E45').Value = '  +6.48%  '
$ws.Range('D46').Value = '''62.08'
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('D47').Value = '''4.87'
$ws.Range('E47').Value = '  -4.10%  '
$ws.Range('E48').Value = '  +3.88%  '
$ws.Range('B49').Value = 'BinanceUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D49').Value = '''1.00'
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').Value = '''1.18'
$ws.Range('E50').Value = '  +3.10%  '
$ws.Range('D51').Value = '''99.13'
$ws.Range('E51').Value = '  +6.50%  '
